############################################################
# "speed up loadcases generators by disabling UI"
#
# Actions sheet: updated DesignFactor values for TLO Traffic,
# Wind and Thermal variable actions; VariableLoadcases /
# PermanentLoadcases gain extra loadcase-name columns to cope
# with the new Thermal Uniform/Gradient, Breaking And Traction
# and Earth pressure actions (replacing the old single "Thermal
# Envelope" / combined loadcases).
############################################################

$wb = $excel.ActiveWorkbook

######################################################################
# 1. Actions sheet - revised partial-factors for the variable actions
######################################################################
$wsActions = $wb.Worksheets.Item("Actions")

# TLO Traffic design factor 1.5 -> 1.35
$wsActions.Range("G3").Value = 1.35
# Wind design factor 1.5 -> 1.75
$wsActions.Range("G5").Value = 1.75
# Thermal design factor 1.5 -> 1.55
$wsActions.Range("G6").Value = 1.55

######################################################################
# 2. VariableLoadcases - split "Thermal Envelope" into Thermal Uniform
#    / Thermal Gradient, and add Breaking And Traction / Horizontal
#    Envelope loadcases to the "TLO Traffic" row
######################################################################
$wsVar = $wb.Worksheets.Item("VariableLoadcases")

$wsVar.Range("B5").Value = "Thermal Uniform"
$wsVar.Range("C5").Value = "Thermal Gradient"

$wsVar.Range("C2").Value = "Breaking And Traction"

######################################################################
# 3. PermanentLoadcases - add Earth pressure (horizontal/vertical)
#    loadcases to the "Permanent load" row
######################################################################
$wsPerm = $wb.Worksheets.Item("PermanentLoadcases")

# extend the header row formatting into the new columns
$wsPerm.Range("C1").Copy($wsPerm.Range("D1:E1"))
$wsPerm.Range("D1").Value = "Loadcase Name"
$wsPerm.Range("E1").Value = "Loadcase Name"

$wsPerm.Range("D2").Value = "Earth pressure horizontal"
$wsPerm.Range("E2").Value = "Earth pressure vertical"

# widen the columns to fit the new, longer loadcase names
$wsPerm.Columns.Item(3).ColumnWidth = 18
$wsPerm.Columns.Item(4).ColumnWidth = 23.166666666666668
$wsPerm.Columns.Item(5).ColumnWidth = 23.166666666666668

# back to VariableLoadcases for the remaining new loadcase and column sizing
$wsVar.Range("D2").Value = "Horizontal Envelope"

# unify the loadcase-name columns to a single, wider width
$wsVar.Range("B1:D1").ColumnWidth = 20.166666666666668

######################################################################
# 4. Selections / active sheet - UI state left by the author after
#    editing (VariableLoadcases ends up the active tab)
######################################################################
$wsActions.Range("G25").Select()
$wsPerm.Range("D13").Select()

$wsAcc = $wb.Worksheets.Item("AccidentalLoadcases")
$wsAcc.Range("C8").Select()

$wsSeismic = $wb.Worksheets.Item("SeismicLoadcases")
$wsSeismic.Range("G12").Select()

$wsVar.Range("I17").Select()
